$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 256.5101652550508
$ws.Range("C2").Value = 95.76379502855229
$ws.Range("D2").Value = 95.76379502855229
$ws.Range("E2").Value = 95.76379502855229
